$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column default formatting (applies to the whole column / future cells) ----
$ws.Columns.Item("A").HorizontalAlignment = -4131
$ws.Columns.Item("B").HorizontalAlignment = -4152
$ws.Columns.Item("C").HorizontalAlignment = -4152
$ws.Columns.Item("C").NumberFormat = "#,##0"
$ws.Columns.Item("D").HorizontalAlignment = -4152
$ws.Columns.Item("D").NumberFormat = "#,##0"
$ws.Columns.Item("E").HorizontalAlignment = 1
$ws.Columns.Item("F").HorizontalAlignment = 1
$ws.Columns.Item("G").HorizontalAlignment = -4152
$ws.Columns.Item("H").HorizontalAlignment = -4152
$ws.Columns.Item("I").HorizontalAlignment = 1
$ws.Columns.Item("J").HorizontalAlignment = -4152
$ws.Columns.Item("J").NumberFormat = "#,##0"
$ws.Columns.Item("K").HorizontalAlignment = 1
$ws.Columns.Item("L").HorizontalAlignment = -4152
$ws.Columns.Item("L").NumberFormat = "#,##0"
$ws.Columns.Item("M").HorizontalAlignment = -4152

# ---- Column widths ----
$ws.Columns.Item("A").ColumnWidth = 16.14785714285714
$ws.Columns.Item("B").ColumnWidth = 14.862142857142858
$ws.Columns.Item("C").ColumnWidth = 16.290714285714284
$ws.Columns.Item("D").ColumnWidth = 14.290714285714287
$ws.Columns.Item("E").ColumnWidth = 18.290714285714284
$ws.Columns.Item("F").ColumnWidth = 23.005
$ws.Columns.Item("G").ColumnWidth = 16.576428571428572
$ws.Columns.Item("H").ColumnWidth = 15.147857142857141
$ws.Columns.Item("I").ColumnWidth = 14.147857142857141
$ws.Columns.Item("J").ColumnWidth = 14.147857142857141
$ws.Columns.Item("K").ColumnWidth = 15.290714285714287
$ws.Columns.Item("L").ColumnWidth = 16.862142857142857
$ws.Columns.Item("M").ColumnWidth = 14.147857142857141

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "email"
$ws.Range("C1").Value = "phoneNumber"
$ws.Range("D1").Value = "paymentDay"
$ws.Range("E1").Value = "collector"
$ws.Range("F1").Value = "level3Address"
$ws.Range("G1").Value = "level2Address"
$ws.Range("H1").Value = "level1Address"
$ws.Range("I1").Value = "plan"
$ws.Range("J1").Value = "planPrice"
$ws.Range("K1").Value = "pricePerCounter"
$ws.Range("L1").Value = "lastCounterValue"
$ws.Range("M1").Value = "isPerCounter"

# ---- Data row (row 2) ----
$ws.Range("A2").Value = "Customer Name"
$ws.Range("B2").Value = "optional_email@email.com"
$ws.Range("C2").Value = 96171234567
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = "employeeUsername"
$ws.Range("F2").Value = "Street"
$ws.Range("G2").Value = "Building"
$ws.Range("H2").Value = "Floor"
$ws.Range("I2").Value = "4MB"
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 0.4
$ws.Range("L2").Value = 2419
$ws.Range("M2").Value = "TRUE/FALSE"

# ---- Per-cell formatting ----
$ws.Range("A1").HorizontalAlignment = -4131
$ws.Range("B1").HorizontalAlignment = -4152
$ws.Range("C1").Font.Name = "Calibri"
$ws.Range("C1").Font.Color = 0
$ws.Range("C1").HorizontalAlignment = -4152
$ws.Range("C1").NumberFormat = "#,##0"
$ws.Range("D1").Font.Name = "Calibri"
$ws.Range("D1").Font.Color = 0
$ws.Range("D1").HorizontalAlignment = -4152
$ws.Range("D1").NumberFormat = "#,##0"
$ws.Range("E1").HorizontalAlignment = 1
$ws.Range("F1").HorizontalAlignment = 1
$ws.Range("G1").HorizontalAlignment = -4152
$ws.Range("H1").Font.Name = "Calibri"
$ws.Range("H1").Font.Color = 0
$ws.Range("H1").HorizontalAlignment = -4152
$ws.Range("I1").HorizontalAlignment = 1
$ws.Range("J1").Font.Name = "Calibri"
$ws.Range("J1").Font.Color = 0
$ws.Range("J1").HorizontalAlignment = -4131
$ws.Range("J1").NumberFormat = "#,##0"
$ws.Range("K1").Font.Name = "Calibri"
$ws.Range("K1").Font.Color = 0
$ws.Range("K1").HorizontalAlignment = -4131
$ws.Range("L1").HorizontalAlignment = -4152
$ws.Range("L1").NumberFormat = "#,##0"
$ws.Range("M1").Font.Name = "Calibri"
$ws.Range("M1").Font.Color = 0
$ws.Range("M1").Font.Underline = $true
$ws.Range("M1").HorizontalAlignment = -4131
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.Color = 0
$ws.Range("A2").HorizontalAlignment = -4152
$ws.Range("B2").Font.Name = "Calibri"
$ws.Range("B2").Font.Color = 0
$ws.Range("B2").HorizontalAlignment = -4152
$ws.Range("B2").NumberFormat = "#,##0"
$ws.Range("C2").Font.Name = "Calibri"
$ws.Range("C2").Font.Color = 0
$ws.Range("C2").HorizontalAlignment = -4152
$ws.Range("C2").NumberFormat = "#,##0"
$ws.Range("D2").HorizontalAlignment = -4152
$ws.Range("D2").NumberFormat = "#,##0"
$ws.Range("E2").HorizontalAlignment = 1
$ws.Range("F2").Font.Name = "Calibri"
$ws.Range("F2").Font.Color = 0
$ws.Range("F2").HorizontalAlignment = -4152
$ws.Range("G2").Font.Name = "Calibri"
$ws.Range("G2").Font.Color = 0
$ws.Range("G2").HorizontalAlignment = -4152
$ws.Range("H2").HorizontalAlignment = -4152
$ws.Range("I2").Font.Name = "Calibri"
$ws.Range("I2").Font.Color = 0
$ws.Range("I2").HorizontalAlignment = -4152
$ws.Range("I2").NumberFormat = "#,##0"
$ws.Range("J2").Font.Name = "Calibri"
$ws.Range("J2").Font.Color = 0
$ws.Range("J2").HorizontalAlignment = -4152
$ws.Range("J2").NumberFormat = "#,##0"
$ws.Range("K2").Font.Name = "Calibri"
$ws.Range("K2").Font.Color = 0
$ws.Range("K2").HorizontalAlignment = -4152
$ws.Range("K2").NumberFormat = "#,##0.00"
$ws.Range("L2").HorizontalAlignment = -4152
$ws.Range("L2").NumberFormat = "#,##0"
$ws.Range("M2").HorizontalAlignment = -4152
